# Enterprise Hierarchy location list maintenance:
#  - remove "Colombo 05" (row 3)
#  - add three new locations at the bottom: Port City Colombo, Negombo, Thalahena
#  - bold the header ("Location")
#  - set print orientation to Portrait
#  - leave the selection where the user ended up after the edits

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Colombo 05" row - everything below shifts up automatically.
$ws.Rows.Item(3).Delete() | Out-Null

# Append the new locations under "Colombo City Centre (CCC)" (now row 9).
$ws.Range("A10").Value = "Port City Colombo"
$ws.Range("A11").Value = "Negombo"
$ws.Range("A12").Value = "Thalahena"

# Make the header stand out.
$ws.Range("A1").Font.Bold = $true

# Print setup -> portrait.
$ws.PageSetup.Orientation = 1

# Match the final on-screen selection.
$ws.Range("A15").Select() | Out-Null
